$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5829956567556286
$ws.Range("D2").Value = 0.565824936022151

$ws.Range("C3").Value = 0.007261972822570133
$ws.Range("D3").Value = 0.9942712849491435

$ws.Range("C4").Value = 0.9236304513466169
$ws.Range("D4").Value = 0.3656988547037152

$ws.Range("C5").Value = -0.1557742007531301
$ws.Range("D5").Value = 0.8776315125607972

$ws.Range("C6").Value = -0.3225531089082612
$ws.Range("D6").Value = 0.7500790726291151

$ws.Range("C7").Value = 0.3231987371382763
$ws.Range("D7").Value = 0.7495967044289018

$ws.Range("C8").Value = -0.5649576966487505
$ws.Range("D8").Value = 0.5778186555861351

$ws.Range("C9").Value = 0.6901467025507328
$ws.Range("D9").Value = 0.4973186532259661

$ws.Range("C10").Value = -0.08135500577877723
$ws.Range("D10").Value = 0.935895236898034

$ws.Range("C11").Value = -0.9538440776757816
$ws.Range("D11").Value = 0.350526343944447
